# Modified Original Course Code
# Adds four new rows to the "GIT" worksheet describing git fork/clone/pull/push
# semantics, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GIT")

# Populate the new cells in the same order the shared-string table grows in
# the target workbook (new unique strings are appended as they are first
# encountered): "git fork", then the three explanatory sentences, in the
# order B19, B21, B18, B20 - followed by the reused command labels.
$ws.Range("A19").Value = "git fork"
$ws.Range("B19").Value = "means you are copying the repository to your Github account"
$ws.Range("B21").Value = "means you are returning the repository after modifying it"
$ws.Range("B18").Value = "means you are making a copy of the repository in your system == DOWNLOADING"
$ws.Range("B20").Value = "means you are fetching the last modified repository == REFRESHING"

$ws.Range("A18").Value = "git clone"
$ws.Range("A20").Value = "git pull"
$ws.Range("A21").Value = "git push"

# Match the row heights present in the target worksheet.
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 28.8
$ws.Rows.Item(20).RowHeight = 28.8
$ws.Rows.Item(21).RowHeight = 28.8

# Move the selection to the last populated cell, as in the target view.
$ws.Range("B21").Select()
